$d = $word.ActiveDocument

# Replace the ID placeholder text (collapses the two runs -
# "**ID__AFFARS_5328_topic_6__ID**" + a trailing space run - into one run
# with the new ID text).
$d.Content.Find.Execute("**ID__AFFARS_5328_topic_6__ID** ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "**ID__AFFARS_5328_106_6__ID**", 2)

# Update the first paragraph's formatting: add a paragraph border (space-only,
# no line) and widen the left indent from 120 twips (6pt) to 225 twips (11.25pt).
$p1 = $d.Paragraphs.Item(1)
$p1.Format.LeftIndent = 11.25
$p1.Format.Borders.DistanceFromTop = 5
$p1.Format.Borders.DistanceFromLeft = 5
$p1.Format.Borders.DistanceFromBottom = 5
$p1.Format.Borders.DistanceFromRight = 5
